# Daily attendance processing - 2026-02-01 22:03:53
#
# For every data row in the "Recorded By" column (G), when the cell holds a
# comma-separated list of recorders, swap the first two entries (leaving any
# further entries, e.g. the trailing "system" on multi-recorder rows, in
# place). Cells holding a single recorder (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $parts = $text -split ", "
    if ($parts.Count -ge 2) {
        $first = $parts[0]
        $second = $parts[1]
        $parts[0] = $second
        $parts[1] = $first
        $cell.Value = ($parts -join ", ")
    }
}
